$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the status of project in row 5 from "approved" to "closed"
$ws.Range("F5").Value = "closed"

# 2. Adjust some column widths (B, C lose bestFit/autosize tweak, widen a bit; K widened)
$ws.Columns.Item(2).ColumnWidth = 13
$ws.Columns.Item(3).ColumnWidth = 13.36328125
$ws.Columns.Item(11).ColumnWidth = 24.54296875

# 3. Append three new project rows (17-19)
$ws.Range("A17").Value = "GG2574529"
$ws.Range("B17").Value = 14.600127799999999
$ws.Range("C17").Value = -91.226982599999999
$ws.Range("E17").Value = "Reforesting Santiago"
$ws.Range("F17").Value = "proposed"
$ws.Range("G17").Value = "Mike"
$ws.Range("H17").Value = "This project restores ecosystems, provides families with sustainable source of firewood, imporoves access to clean water and strengthens community resilience."
$ws.Range("I17").Value = 38000
$ws.Range("J17").Value = "Rotary E-Club San Diego Global"
$ws.Range("K17").Value = 5340
$ws.Range("L17").Value = "WKG"
$ws.Range("M17").Value = 2025

$ws.Range("A18").Value = "ADP_VistaHermosa"
$ws.Range("B18").Value = 14.777010000000001
$ws.Range("C18").Value = -91.267499999999998
$ws.Range("D18").Value = "https://docs.google.com/document/d/1piFpvjaRl6BletKqpdrz2wlrc7pAYXgqHn6leAGrNto/edit?tab=t.0"
$ws.Range("E18").Value = "Sanitation for Vista Hermosa, Santa Lucia Utatlan"
$ws.Range("F18").Value = "proposed"
$ws.Range("G18").Value = "Bruce"
$ws.Range("H18").Value = "This project povides sanitation services for 30 familes.  This is a follow-on project to an earlier grant that rehabiliated a potable water system."
$ws.Range("I18").Value = 50000
$ws.Range("L18").Value = "AdP"
$ws.Range("M18").Value = 2025

$ws.Range("A19").Value = "ADP_Panimaquip"
$ws.Range("B19").Value = 14.600020000000001
$ws.Range("C19").Value = -91.147509999999997
$ws.Range("D19").Value = "https://docs.google.com/document/d/1DhmMi5CsALaOjiPloV8raZ8VPry1vaK8zJoh0PWGvdo/edit?tab=t.0"
$ws.Range("E19").Value = "Water and Sanitation for Panimaquip, San Lucas Toliman"
$ws.Range("F19").Value = "proposed"
$ws.Range("G19").Value = "Bruce"
$ws.Range("H19").Value = "This project rehabiltates a potable water system."
$ws.Range("I19").Value = 50000
$ws.Range("L19").Value = "AdP"
$ws.Range("M19").Value = 2025

# 4. Turn the data range into an Excel Table ("Table1")
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:M19"), [System.Type]::Missing, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = "TableStyleLight1"

# 5. Update view state (top-left cell + active selection)
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("E1").Select()
$ws.Range("I30").Select()
